# Convert charts to altair
# The header cell that used to read "Type" now reads "Housing" (the
# companion header cell "Usage" stays the same). Update the cell value
# and move the active selection to D3, matching the authored edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Housing"

$ws.Range("D3").Select()
